# Generate Report for Archive
#
# The localization status "Ready for handoff" has moved on to
# "In Translation" everywhere it is reported: the per-language Overview
# columns (zh-cn / de-de) and the Status column on each language's own
# sheet. Once the text is shorter, the Status-ish columns are re-sized to
# fit the new value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview!E2:F3 -> zh-cn / de-de status for both rows of content
$overview.Range("E2:F3").Value = "In Translation"

# Each language sheet's own Status column (C) for both rows of content
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# Re-fit the now-narrower status columns to the new text
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
